$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 9's date column: change A9 from 21/05/2021 to 17/06/2024
# (this was stored as literal text, not a real date, so keep it text)
$ws.Range("A9").Value = "17/06/2024"

# Leave the selection on A10 (matches final cursor position in the authoring session)
$null = $ws.Range("A10").Select()
